# Alpha 0.3.0 - Teeth, Tanks, and the Red Thief have arrived - all enemies
# are now in the game! Update the CC Feature List sheet: mark the enemy
# tiles as Complete, rename the "Random" enemy to "Walker", flip a couple
# of statuses, and add "(needs testing w/ enemies)" notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enemy column (G/H): every enemy implementation is now Complete.
$ws.Range("H2").Value = "Complete"
$ws.Range("H3").Value = "Complete"
$ws.Range("H4").Value = "Complete"
$ws.Range("H5").Value = "Complete"
$ws.Range("H6").Value = "Complete"
$ws.Range("H7").Value = "Complete"

# "Random" enemy was renamed to "Walker".
$ws.Range("G8").Value = "Walker"
$ws.Range("H8").Value = "Complete"

$ws.Range("H9").Value = "Complete"
$ws.Range("H10").Value = "Complete"

# Gravel is now Complete.
$ws.Range("B14").Value = "Complete"

# Bomb note: needs testing with enemies.
$ws.Range("C16").Value = "(needs testing w/ enemies)"

# Teleport moved from Inactive to Prototype, also needs testing w/ enemies.
$ws.Range("B19").Value = "Prototype"
$ws.Range("C19").Value = "(needs testing w/ enemies)"

# Tank Button: "See Enemies" status item removed, now Complete.
$ws.Range("B21").Value = "Complete"

# Restore the author's last selection.
$ws.Range("E13").Select()
